$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(17, 9).Value = "b"
$ws.Cells.Item(17, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(22, 9).Value = "sv"
$ws.Cells.Item(22, 10).Value = "Statement-opinion"
$ws.Cells.Item(27, 9).Value = "sv"
$ws.Cells.Item(27, 10).Value = "Statement-opinion"
$ws.Cells.Item(28, 9).Value = "sd"
$ws.Cells.Item(28, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(31, 9).Value = "sv"
$ws.Cells.Item(31, 10).Value = "Statement-opinion"
$ws.Cells.Item(42, 9).Value = "sv"
$ws.Cells.Item(42, 10).Value = "Statement-opinion"
$ws.Cells.Item(43, 9).Value = "sd"
$ws.Cells.Item(43, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(56, 9).Value = "aa"
$ws.Cells.Item(56, 10).Value = "Agree/Accept"
$ws.Cells.Item(61, 9).Value = "sv"
$ws.Cells.Item(61, 10).Value = "Statement-opinion"
$ws.Cells.Item(63, 9).Value = "ba"
$ws.Cells.Item(63, 10).Value = "Appreciation"
$ws.Cells.Item(71, 9).Value = "sv"
$ws.Cells.Item(71, 10).Value = "Statement-opinion"
$ws.Cells.Item(77, 9).Value = "qy"
$ws.Cells.Item(77, 10).Value = "Yes-No-Question"
$ws.Cells.Item(82, 9).Value = "sv"
$ws.Cells.Item(82, 10).Value = "Statement-opinion"
$ws.Cells.Item(99, 9).Value = "aa"
$ws.Cells.Item(99, 10).Value = "Agree/Accept"
$ws.Cells.Item(108, 9).Value = "b"
$ws.Cells.Item(108, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(110, 9).Value = "ba"
$ws.Cells.Item(110, 10).Value = "Appreciation"
$ws.Cells.Item(120, 9).Value = "sd"
$ws.Cells.Item(120, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(121, 9).Value = "aa"
$ws.Cells.Item(121, 10).Value = "Agree/Accept"
$ws.Cells.Item(124, 9).Value = "sv"
$ws.Cells.Item(124, 10).Value = "Statement-opinion"
$ws.Cells.Item(134, 9).Value = "aa"
$ws.Cells.Item(134, 10).Value = "Agree/Accept"
$ws.Cells.Item(135, 9).Value = "b"
$ws.Cells.Item(135, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(136, 9).Value = "aa"
$ws.Cells.Item(136, 10).Value = "Agree/Accept"
$ws.Cells.Item(141, 9).Value = "sd"
$ws.Cells.Item(141, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(146, 9).Value = "sd"
$ws.Cells.Item(146, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(149, 9).Value = "b"
$ws.Cells.Item(149, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(151, 9).Value = "sv"
$ws.Cells.Item(151, 10).Value = "Statement-opinion"
$ws.Cells.Item(152, 9).Value = "sd"
$ws.Cells.Item(152, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(173, 9).Value = "aa"
$ws.Cells.Item(173, 10).Value = "Agree/Accept"
$ws.Cells.Item(200, 9).Value = "b"
$ws.Cells.Item(200, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(211, 9).Value = "sd"
$ws.Cells.Item(211, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(214, 9).Value = "ba"
$ws.Cells.Item(214, 10).Value = "Appreciation"
$ws.Cells.Item(220, 9).Value = "sd"
$ws.Cells.Item(220, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(230, 9).Value = "aa"
$ws.Cells.Item(230, 10).Value = "Agree/Accept"
$ws.Cells.Item(233, 9).Value = "sd"
$ws.Cells.Item(233, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(246, 9).Value = "sd"
$ws.Cells.Item(246, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(248, 9).Value = "sv"
$ws.Cells.Item(248, 10).Value = "Statement-opinion"
$ws.Cells.Item(256, 9).Value = "sd"
$ws.Cells.Item(256, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(267, 9).Value = "%"
$ws.Cells.Item(267, 10).Value = "Uninterpretable"
$ws.Cells.Item(272, 9).Value = "sd"
$ws.Cells.Item(272, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(288, 9).Value = "aa"
$ws.Cells.Item(288, 10).Value = "Agree/Accept"
$ws.Cells.Item(291, 9).Value = "%"
$ws.Cells.Item(291, 10).Value = "Uninterpretable"
$ws.Cells.Item(315, 9).Value = "sd"
$ws.Cells.Item(315, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(318, 9).Value = "ba"
$ws.Cells.Item(318, 10).Value = "Appreciation"
$ws.Cells.Item(358, 9).Value = "sd"
$ws.Cells.Item(358, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(360, 9).Value = "sv"
$ws.Cells.Item(360, 10).Value = "Statement-opinion"
$ws.Cells.Item(362, 9).Value = "sd"
$ws.Cells.Item(362, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(367, 9).Value = "b"
$ws.Cells.Item(367, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(372, 9).Value = "b"
$ws.Cells.Item(372, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(374, 9).Value = "b"
$ws.Cells.Item(374, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(378, 9).Value = "sv"
$ws.Cells.Item(378, 10).Value = "Statement-opinion"
$ws.Cells.Item(446, 9).Value = "aa"
$ws.Cells.Item(446, 10).Value = "Agree/Accept"
$ws.Cells.Item(457, 9).Value = "sv"
$ws.Cells.Item(457, 10).Value = "Statement-opinion"
$ws.Cells.Item(459, 9).Value = "sv"
$ws.Cells.Item(459, 10).Value = "Statement-opinion"
$ws.Cells.Item(462, 9).Value = "sv"
$ws.Cells.Item(462, 10).Value = "Statement-opinion"
$ws.Cells.Item(485, 9).Value = "aa"
$ws.Cells.Item(485, 10).Value = "Agree/Accept"
$ws.Cells.Item(489, 9).Value = "sd"
$ws.Cells.Item(489, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(492, 9).Value = "aa"
$ws.Cells.Item(492, 10).Value = "Agree/Accept"
$ws.Cells.Item(498, 9).Value = "b"
$ws.Cells.Item(498, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(502, 9).Value = "ba"
$ws.Cells.Item(502, 10).Value = "Appreciation"
$ws.Cells.Item(504, 9).Value = "aa"
$ws.Cells.Item(504, 10).Value = "Agree/Accept"
$ws.Cells.Item(506, 9).Value = "ba"
$ws.Cells.Item(506, 10).Value = "Appreciation"
$ws.Cells.Item(508, 9).Value = "sv"
$ws.Cells.Item(508, 10).Value = "Statement-opinion"
$ws.Cells.Item(516, 9).Value = "sd"
$ws.Cells.Item(516, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(518, 9).Value = "ba"
$ws.Cells.Item(518, 10).Value = "Appreciation"
$ws.Cells.Item(522, 9).Value = "b"
$ws.Cells.Item(522, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(523, 9).Value = "sd"
$ws.Cells.Item(523, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(528, 9).Value = "sv"
$ws.Cells.Item(528, 10).Value = "Statement-opinion"
$ws.Cells.Item(531, 9).Value = "sd"
$ws.Cells.Item(531, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(532, 9).Value = "sv"
$ws.Cells.Item(532, 10).Value = "Statement-opinion"
$ws.Cells.Item(537, 9).Value = "sv"
$ws.Cells.Item(537, 10).Value = "Statement-opinion"
$ws.Cells.Item(546, 9).Value = "sv"
$ws.Cells.Item(546, 10).Value = "Statement-opinion"
$ws.Cells.Item(558, 9).Value = "sv"
$ws.Cells.Item(558, 10).Value = "Statement-opinion"
$ws.Cells.Item(564, 9).Value = "ba"
$ws.Cells.Item(564, 10).Value = "Appreciation"
$ws.Cells.Item(566, 9).Value = "ba"
$ws.Cells.Item(566, 10).Value = "Appreciation"
$ws.Cells.Item(568, 9).Value = "sv"
$ws.Cells.Item(568, 10).Value = "Statement-opinion"
$ws.Cells.Item(569, 9).Value = "ba"
$ws.Cells.Item(569, 10).Value = "Appreciation"
Write-Output "Updated dialog act annotations for 76 rows."
